$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28: [GNN 2] Graph LSTM -> Mobile manipulator 101 :: Basic Kinematics
$ws.Range("D28").Value = "Mobile manipulator 101 :: Basic Kinematics"
$ws.Range("E28").Value = "https://ropiens.tistory.com/196"

# Row 32: Graph 유형 정리 (GNN) -> Hadoop Ecosystem 하둡 에코시스템 간단정리
$ws.Range("D32").Value = "Hadoop Ecosystem 하둡 에코시스템 간단정리"
$ws.Range("E32").Value = "https://dodonam.tistory.com/390"

# Row 51: [pandas] 데이터프레임 행 또는 컬럼 삭제하기, drop() 메소드 -> [pandas] 컬럼 값이 특정 조건에 부합하는 행들의 특정 컬럼 값 수정하기
$ws.Range("D51").Value = "[pandas] 컬럼 값이 특정 조건에 부합하는 행들의 특정 컬럼 값 수정하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-%EC%BB%AC%EB%9F%BC-%EA%B0%92%EC%9D%B4-%ED%8A%B9%EC%A0%95-%EC%A1%B0%EA%B1%B4%EC%97%90-%EB%B6%80%ED%95%A9%ED%95%98%EB%8A%94-%ED%96%89%EB%93%A4%EC%9D%98-%ED%8A%B9%EC%A0%95-%EC%BB%AC%EB%9F%BC-%EA%B0%92-%EC%88%98%EC%A0%95%ED%95%98%EA%B8%B0"
